$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Amount Due" column (old column H) entirely - this shifts
# the old "Serial Number" column (I) left into H, and shrinks the used
# range from A1:I1000 to A1:H1000.
$ws.Columns("H:H").Delete()

# Re-point the headers (row 1) to the new column order / new wording.
$ws.Range("A1").Value2 = "Booking ID"
$ws.Range("B1").Value2 = "Appliance"
$ws.Range("C1").Value2 = "Spare Part"
$ws.Range("D1").Value2 = "Parts Number"
$ws.Range("E1").Value2 = "Quantity"
$ws.Range("F1").Value2 = "Age of Requested"
$ws.Range("G1").Value2 = "Model No"
$ws.Range("H1").Value2 = "Serial No"

# Re-point the merge-field placeholders (row 2) to match the new column
# order so the downloaded file's columns follow the same sequence as the
# on-screen view.
$ws.Range("A2").Value2 = "{spare:booking_id}"
$ws.Range("B2").Value2 = "{spare:services}"
$ws.Range("C2").Value2 = "{spare:parts_requested}"
$ws.Range("D2").Value2 = "{spare:part_number}"
$ws.Range("E2").Value2 = "{spare:quantity}"
$ws.Range("F2").Value2 = "{spare:age_of_request}"
$ws.Range("G2").Value2 = "{spare:model_number}"
$ws.Range("H2").Value2 = "{spare:serial_number}"

# Put the selection/view back at the top-left (H1), matching the saved
# view state.
$ws.Range("H1").Select()
